# Adding AccountPage data: 8 new user rows (rows 76-83) to the userData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: FirstName, LastName, Email, Password, DateOfBirth
$rows = @(
    @('Albertine', 'Gorczany',  'marlon.koepp@example.com',     '4a6tK!2&klW',  '01/27/1981'),
    @('Lacresha',  'Johnston',  'youlanda.carter@example.com',  '^$iI7I',       '07/24/1991'),
    @('Glenn',     'Zemlak',    'loris.padberg@example.com',    'JRF4094AF',    '09/04/1956'),
    @('Johnnie',   'Toy',       'deloise.hammes@example.com',   '4Ig7NHH4c',    '05/29/1986'),
    @('Paris',     'Mills',     'jamison.olson@example.com',    '&*@3X9@d',     '01/05/1967'),
    @('Brice',     'Will',      'emmie.fahey@example.com',      'zsRV#sjw',     '10/04/1982'),
    @('Beaulah',   'Lehner',    'ka.corkery@example.com',       '7Tb4139KuA',   '06/11/1988'),
    @('Peter',     'Reichert',  'garry.mccullough@example.com', '2#v@3x1',      '01/22/1964')
)

$startRow = 76
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]

    # The last column holds a MM/DD/YYYY-looking string that must stay plain
    # text (like the rest of the sheet) rather than being auto-converted to
    # a date serial number. Entering it as a formula that evaluates to the
    # literal string, then collapsing the formula to its value via
    # copy/paste-special, keeps the cell a genuine text value without
    # touching any cell's number format/style.
    $dateCell = $ws.Cells.Item($r, 5)
    $dateCell.Formula = '="' + $data[4] + '"'
    $dateCell.Copy($dateCell) | Out-Null
    $dateCell.PasteSpecial(-4163) | Out-Null
}
